$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$wsQ2Name = $wb.Worksheets.Item(2).Name

# A truly blank, default-styled cell used below to strip formatting back
# to the workbook default (style index 0) after forcing a numeric-looking
# string to stay text.
$blank = $ws1.Range("Z100")

function Set-TextValue($range, $value) {
    # Force the assignment to be stored as text (otherwise Excel silently
    # coerces numeric-looking strings like "002379" or "1.58" into numbers
    # and drops leading zeros / trailing zeros), then strip the number
    # format back off so the cell keeps the sheet's default style.
    $range.NumberFormat = "@"
    $range.Value = $value
    $blank.Copy()
    $range.PasteSpecial(-4122)  # xlPasteFormats
}

# ------------------------------------------------------------------
# 1) Update sheet "总计": row 2 becomes the new 2022-Q3 total, and a
#    new row 3 is appended carrying the former 2022-Q2 total down.
# ------------------------------------------------------------------

# Duplicate row 2's formatting into row 3 first (keeps identical styles),
# then overwrite the values that actually change.
$ws1.Range("A2:D2").Copy($ws1.Range("A3:D3"))

$ws1.Range("A3").Value = 1
$ws1.Range("B3").Value = "2022-Q2"
$ws1.Range("C3").Value = 2
$ws1.Range("D3").Value = 0.14

$ws1.Range("B2").Value = "2022-Q3"
$ws1.Range("D2").Value = 0.17

# ------------------------------------------------------------------
# 2) Insert a brand-new worksheet "2022-Q3" between "总计" and the
#    existing "2022-Q2" sheet, holding the quarterly fund detail.
# ------------------------------------------------------------------
$wsQ2 = $wb.Worksheets.Item($wsQ2Name)
$wsQ3 = $wb.Worksheets.Add($wsQ2, $null)
$wsQ3.Name = "2022-Q3"

# Header row - copy style from sheet1's header row so formatting matches.
$ws1.Range("B1:D1").Copy($wsQ3.Range("B1:D1"))
$ws1.Range("B1:D1").Copy($wsQ3.Range("E1:G1"))
$ws1.Range("B1").Copy($wsQ3.Range("H1"))

$wsQ3.Range("B1").Value = "基金代码"
$wsQ3.Range("C1").Value = "基金名称"
$wsQ3.Range("D1").Value = "基金规模"
$wsQ3.Range("E1").Value = "股票总仓位"
$wsQ3.Range("F1").Value = "仓位占比"
$wsQ3.Range("G1").Value = "持有市值(亿元)"
$wsQ3.Range("H1").Value = "仓位排名"

# "A" column cells - copy style from sheet1's A2 (same look as s="2").
$ws1.Range("A2").Copy($wsQ3.Range("A2"))
$ws1.Range("A2").Copy($wsQ3.Range("A3"))
$wsQ3.Range("A2").Value = 0
$wsQ3.Range("A3").Value = 1

# Row 2 - numeric-looking text columns must remain text.
Set-TextValue $wsQ3.Range("B2") "002379"
Set-TextValue $wsQ3.Range("C2") "工银瑞信香港中小盘股票（QDII）人民币"
Set-TextValue $wsQ3.Range("D2") "1.58"
Set-TextValue $wsQ3.Range("E2") "78.58"
Set-TextValue $wsQ3.Range("F2") "5.44"
Set-TextValue $wsQ3.Range("G2") "0.0860"
$wsQ3.Range("H2").Value = 3

# Row 3
Set-TextValue $wsQ3.Range("B3") "002380"
Set-TextValue $wsQ3.Range("C3") "工银瑞信香港中小盘股票（QDII）美元"
Set-TextValue $wsQ3.Range("D3") "1.58"
Set-TextValue $wsQ3.Range("E3") "78.58"
Set-TextValue $wsQ3.Range("F3") "5.44"
Set-TextValue $wsQ3.Range("G3") "0.0860"
$wsQ3.Range("H3").Value = 3

# Keep the originally-active "2022-Q2" sheet as the selected tab, same as
# before the edit (it merely shifted from position 2 to position 3).
$wb.Worksheets.Item($wsQ2Name).Activate()

Write-Host "done"
